# Slide 9 ("What is git?") figure update.
#
# - Shrinks the screenshot picture.
# - Re-purposes the existing red arrow + red caption (which pointed at the
#   "commit" step) to point at a "push" step lower on the screenshot, with
#   new wording ("Push modified repo to Github").
# - Adds a second red arrow + caption ("Clone repo from Github") higher up,
#   duplicated from the first arrow/caption pair.

# Helper: PowerPoint's Shape.Left/Top/Width/Height are single-precision
# (float32) point values. Converting target EMUs -> points -> (float32)
# and back can be off by 1 EMU after truncation; nudge the point value by
# the smallest amount needed so the round trip lands exactly on the target
# EMU value.
function Emu2Pt([double]$emu) {
    $base = $emu / 12700.0
    for ($i = 0; $i -le 200; $i++) {
        $candidate = $base + ($i * 0.0000001)
        $asFloat = [float]$candidate
        $recovered = [math]::Floor(([double]$asFloat * 12700.0) + 0.0000001)
        if ($recovered -eq $emu) {
            return $candidate
        }
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# --- 1. Shrink the "git_crash_course" screenshot picture -------------------
$pic = $s.Shapes.Item(2)
$pic.Width = Emu2Pt 6766560
$pic.Height = Emu2Pt 5178642

# --- 2. Re-aim the existing arrow (id 10) up to the "push" line ------------
$arrow1 = $s.Shapes.Item(3)
$arrow1.HorizontalFlip = $false
$arrow1.VerticalFlip = $true
$arrow1.Left = Emu2Pt 6084805
$arrow1.Top = Emu2Pt 5211227
$arrow1.Width = Emu2Pt 999748
$arrow1.Height = Emu2Pt 5575

# --- 3. Re-word + reposition the existing caption (id 13) ------------------
$caption1 = $s.Shapes.Item(4)
$capRange = $caption1.TextFrame.TextRange
# Replace the first run's text (was "Commit to local repository
# (additional push command needed to upload to ") and drop the trailing
# ")" run so the text reads "Push modified repo to Github".
$capRange.Characters(1, 72).Text = "Push modified repo to "
$capRange.Characters($capRange.Length, 1).Delete()

$caption1.Left = Emu2Pt 7201316
$caption1.Top = Emu2Pt 4888061
$caption1.Width = Emu2Pt 1193180
$caption1.Height = Emu2Pt 1200329

# --- 4. Make PowerPoint hand out shape ids 7/8 for the new duplicates ------
# (PowerPoint reuses the lowest free id on the slide; burn ids 2,3,6 first
# with a few scratch shapes so the upcoming duplicates land on 7 and 8,
# matching the authored deck.)
$scratch = @()
for ($i = 0; $i -lt 3; $i++) {
    $scratch += $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
}
foreach ($j in $scratch) { $j.Delete() }

# --- 5. Duplicate the arrow for the new "clone" callout --------------------
$arrow2 = $arrow1.Duplicate().Item(1)
$arrow2.Name = "Straight Arrow Connector 6"
$arrow2.HorizontalFlip = $true
$arrow2.VerticalFlip = $true
$arrow2.Left = Emu2Pt 6084805
$arrow2.Top = Emu2Pt 3743057
$arrow2.Width = Emu2Pt 999748
$arrow2.Height = Emu2Pt 5575

# --- 6. Duplicate the caption for the new "clone" callout ------------------
$caption2 = $caption1.Duplicate().Item(1)
$caption2.Name = "TextBox 7"
$caption2.Left = Emu2Pt 7201316
$caption2.Top = Emu2Pt 3419891
$caption2.Width = Emu2Pt 1193180
$caption2.Height = Emu2Pt 923330

$cap2Range = $caption2.TextFrame.TextRange
# caption2 currently reads "Push modified repo to Github" (inherited from
# caption1); swap the leading run for the new wording, keeping "Github".
$cap2Range.Characters(1, 22).Text = "Clone repo from "
